$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the product row from the "4" series to the "6" series
$ws.Range("A2").Value = "Test - Dummy Product 6"
$ws.Range("B2").Value = "Dummy Product 6"
$ws.Range("C2").Value = "TEST - Dummy 06"

# Update the saved view: scroll back to column A and select A3
$ws.Range("A3").Select()
